# Update bitcoin_buys.xlsx after running on 2025-10-19
# Appends a new purchase record as row 53 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the purchase date as literal text (matching the existing
# rows below the header, which are plain date strings rather than Excel
# date serials). Force text formatting before assignment so Excel doesn't
# auto-convert the string into a date serial number, then reset the style
# back to the sheet default so no stray formatting is left behind.
$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value = "10/19/2025"
$ws.Range("A53").Style = "Normal"

$ws.Range("B53").Value = 0.0004650000000000036
$ws.Range("C53").Value = 107526.8817204293
$ws.Range("D53").Value = 50
